# Fruta / hortaliza, semanal
# Insert a new weekly price-record row for "Ciboulette" at Mercado Mayorista
# Lo Valledor de Santiago, shifting the existing rows 422..457 down to
# 423..458.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 422 - pushes every row from 422 downward by one.
$ws.Rows.Item(422).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(422, 1).Value = 6
$ws.Cells.Item(422, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(422, 3).Value = "Metropolitana"
$ws.Cells.Item(422, 4).Value = 44746
$ws.Cells.Item(422, 5).Value = 13
$ws.Cells.Item(422, 6).Value = 100112039
$ws.Cells.Item(422, 7).Value = "Ciboulette"
$ws.Cells.Item(422, 8).Value = "Sin especificar"
$ws.Cells.Item(422, 9).Value = "Segunda"
$ws.Cells.Item(422, 10).Value = 280
$ws.Cells.Item(422, 11).Value = 1400
$ws.Cells.Item(422, 12).Value = 1500
$ws.Cells.Item(422, 13).Value = 1443
$ws.Cells.Item(422, 14).Value = "`$/docena de atados"
$ws.Cells.Item(422, 15).Value = "Región Metropolitana"
$ws.Cells.Item(422, 16).Value = 481
$ws.Cells.Item(422, 17).Value = 3
$ws.Cells.Item(422, 18).Value = "Hortaliza"
